$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value2 = 36.86999893188477
$ws.Cells.Item(2, 5).Value2 = 33.40999984741211
$ws.Cells.Item(2, 6).Value2 = 36.97999954223633
$ws.Cells.Item(2, 7).Value2 = 31.14999961853028
$ws.Cells.Item(2, 8).Value2 = 119792704
$ws.Cells.Item(2, 9).Value2 = "PTC"

$ws.Cells.Item(3, 4).Value2 = 36.13000106811523
$ws.Cells.Item(3, 5).Value2 = 38.34000015258789
$ws.Cells.Item(3, 6).Value2 = 38.84999847412109
$ws.Cells.Item(3, 7).Value2 = 35.72000122070312
$ws.Cells.Item(3, 8).Value2 = 119792704
$ws.Cells.Item(3, 9).Value2 = "PTC"

$ws.Cells.Item(4, 4).Value2 = 41.27000045776367
$ws.Cells.Item(4, 5).Value2 = 36.34999847412109
$ws.Cells.Item(4, 6).Value2 = 41.79999923706055
$ws.Cells.Item(4, 7).Value2 = 36
$ws.Cells.Item(4, 8).Value2 = 119792704
$ws.Cells.Item(4, 9).Value2 = "PTC"

$ws.Cells.Item(5, 4).Value2 = 31.76000022888184
$ws.Cells.Item(5, 5).Value2 = 35.43999862670898
$ws.Cells.Item(5, 6).Value2 = 35.70999908447266
$ws.Cells.Item(5, 7).Value2 = 30.53000068664551
$ws.Cells.Item(5, 8).Value2 = 119792704
$ws.Cells.Item(5, 9).Value2 = "PTC"

$ws.Cells.Item(6, 4).Value2 = 34.20000076293945
$ws.Cells.Item(6, 5).Value2 = 29.61000061035156
$ws.Cells.Item(6, 6).Value2 = 34.20000076293945
$ws.Cells.Item(6, 7).Value2 = 28.75
$ws.Cells.Item(6, 8).Value2 = 119792704
$ws.Cells.Item(6, 9).Value2 = "PTC"

$ws.Cells.Item(7, 4).Value2 = 32.84999847412109
$ws.Cells.Item(7, 5).Value2 = 36.45999908447266
$ws.Cells.Item(7, 6).Value2 = 38.5
$ws.Cells.Item(7, 7).Value2 = 31.57999992370605
$ws.Cells.Item(7, 8).Value2 = 119792704
$ws.Cells.Item(7, 9).Value2 = "PTC"

$ws.Cells.Item(8, 4).Value2 = 37.58000183105469
$ws.Cells.Item(8, 5).Value2 = 39.72999954223633
$ws.Cells.Item(8, 6).Value2 = 41.40999984741211
$ws.Cells.Item(8, 7).Value2 = 36.56999969482422
$ws.Cells.Item(8, 8).Value2 = 119792704
$ws.Cells.Item(8, 9).Value2 = "PTC"

$ws.Cells.Item(9, 4).Value2 = 45.02000045776367
$ws.Cells.Item(9, 5).Value2 = 47.43999862670898
$ws.Cells.Item(9, 6).Value2 = 47.63999938964844
$ws.Cells.Item(9, 7).Value2 = 43.09999847412109
$ws.Cells.Item(9, 8).Value2 = 119792704
$ws.Cells.Item(9, 9).Value2 = "PTC"

$ws.Cells.Item(10, 4).Value2 = 46.59000015258789
$ws.Cells.Item(10, 5).Value2 = 52.56999969482422
$ws.Cells.Item(10, 6).Value2 = 54.13000106811523
$ws.Cells.Item(10, 7).Value2 = 45.93000030517578
$ws.Cells.Item(10, 8).Value2 = 119792704
$ws.Cells.Item(10, 9).Value2 = "PTC"

$ws.Cells.Item(11, 4).Value2 = 52.79000091552734
$ws.Cells.Item(11, 5).Value2 = 54.04999923706055
$ws.Cells.Item(11, 6).Value2 = 54.31999969482422
$ws.Cells.Item(11, 7).Value2 = 51
$ws.Cells.Item(11, 8).Value2 = 119792704
$ws.Cells.Item(11, 9).Value2 = "PTC"

$ws.Cells.Item(12, 4).Value2 = 55.40000152587891
$ws.Cells.Item(12, 5).Value2 = 55.18999862670898
$ws.Cells.Item(12, 6).Value2 = 59.29000091552734
$ws.Cells.Item(12, 7).Value2 = 53.45999908447266
$ws.Cells.Item(12, 8).Value2 = 119792704
$ws.Cells.Item(12, 9).Value2 = "PTC"

$ws.Cells.Item(13, 4).Value2 = 56.36999893188477
$ws.Cells.Item(13, 5).Value2 = 66.44999694824219
$ws.Cells.Item(13, 6).Value2 = 66.94999694824219
$ws.Cells.Item(13, 7).Value2 = 55.7400016784668
$ws.Cells.Item(13, 8).Value2 = 119792704
$ws.Cells.Item(13, 9).Value2 = "PTC"

$ws.Cells.Item(14, 4).Value2 = 60.90999984741211
$ws.Cells.Item(14, 5).Value2 = 72.68000030517578
$ws.Cells.Item(14, 6).Value2 = 74.87999725341797
$ws.Cells.Item(14, 7).Value2 = 60.45000076293945
$ws.Cells.Item(14, 8).Value2 = 119792704
$ws.Cells.Item(14, 9).Value2 = "PTC"

$ws.Cells.Item(15, 4).Value2 = 77.62000274658203
$ws.Cells.Item(15, 5).Value2 = 82.34999847412109
$ws.Cells.Item(15, 6).Value2 = 86.48999786376953
$ws.Cells.Item(15, 7).Value2 = 74.76999664306641
$ws.Cells.Item(15, 8).Value2 = 119792704
$ws.Cells.Item(15, 9).Value2 = "PTC"

$ws.Cells.Item(16, 4).Value2 = 93.9499969482422
$ws.Cells.Item(16, 5).Value2 = 91.91000366210938
$ws.Cells.Item(16, 6).Value2 = 99.87000274658205
$ws.Cells.Item(16, 7).Value2 = 90.45999908447266
$ws.Cells.Item(16, 8).Value2 = 119792704
$ws.Cells.Item(16, 9).Value2 = "PTC"

$ws.Cells.Item(17, 4).Value2 = 106.9800033569336
$ws.Cells.Item(17, 5).Value2 = 82.41000366210938
$ws.Cells.Item(17, 6).Value2 = 107.4400024414062
$ws.Cells.Item(17, 7).Value2 = 77.33000183105469
$ws.Cells.Item(17, 8).Value2 = 119792704
$ws.Cells.Item(17, 9).Value2 = "PTC"

$ws.Cells.Item(18, 4).Value2 = 80.87999725341797
$ws.Cells.Item(18, 5).Value2 = 84.79000091552734
$ws.Cells.Item(18, 6).Value2 = 89.33000183105469
$ws.Cells.Item(18, 7).Value2 = 77.26000213623047
$ws.Cells.Item(18, 8).Value2 = 119792704
$ws.Cells.Item(18, 9).Value2 = "PTC"

$ws.Cells.Item(19, 4).Value2 = 93.1999969482422
$ws.Cells.Item(19, 5).Value2 = 90.47000122070312
$ws.Cells.Item(19, 6).Value2 = 102.4700012207031
$ws.Cells.Item(19, 7).Value2 = 89.76999664306641
$ws.Cells.Item(19, 8).Value2 = 119792704
$ws.Cells.Item(19, 9).Value2 = "PTC"

$ws.Cells.Item(20, 4).Value2 = 91.40000152587891
$ws.Cells.Item(20, 5).Value2 = 67.77999877929688
$ws.Cells.Item(20, 6).Value2 = 93.33999633789062
$ws.Cells.Item(20, 7).Value2 = 66.73000335693359
$ws.Cells.Item(20, 8).Value2 = 119792704
$ws.Cells.Item(20, 9).Value2 = "PTC"

$ws.Cells.Item(21, 4).Value2 = 68.69000244140625
$ws.Cells.Item(21, 5).Value2 = 66.91000366210938
$ws.Cells.Item(21, 6).Value2 = 77.75
$ws.Cells.Item(21, 7).Value2 = 62.18999862670898
$ws.Cells.Item(21, 8).Value2 = 119792704
$ws.Cells.Item(21, 9).Value2 = "PTC"

$ws.Cells.Item(22, 4).Value2 = 75.43000030517578
$ws.Cells.Item(22, 5).Value2 = 83.12000274658203
$ws.Cells.Item(22, 6).Value2 = 88.86000061035156
$ws.Cells.Item(22, 7).Value2 = 74.61000061035156
$ws.Cells.Item(22, 8).Value2 = 119792704
$ws.Cells.Item(22, 9).Value2 = "PTC"

$ws.Cells.Item(23, 4).Value2 = 58.61999893188477
$ws.Cells.Item(23, 5).Value2 = 69.25
$ws.Cells.Item(23, 6).Value2 = 71.27999877929688
$ws.Cells.Item(23, 7).Value2 = 52.95999908447266
$ws.Cells.Item(23, 8).Value2 = 119792704
$ws.Cells.Item(23, 9).Value2 = "PTC"

$ws.Cells.Item(24, 4).Value2 = 77.93000030517578
$ws.Cells.Item(24, 5).Value2 = 85.55999755859375
$ws.Cells.Item(24, 6).Value2 = 88.80000305175781
$ws.Cells.Item(24, 7).Value2 = 76.15000152587891
$ws.Cells.Item(24, 8).Value2 = 119792704
$ws.Cells.Item(24, 9).Value2 = "PTC"

$ws.Cells.Item(25, 4).Value2 = 83.66000366210938
$ws.Cells.Item(25, 5).Value2 = 83.87999725341797
$ws.Cells.Item(25, 6).Value2 = 89.44999694824219
$ws.Cells.Item(25, 7).Value2 = 79.36000061035156
$ws.Cells.Item(25, 8).Value2 = 119792704
$ws.Cells.Item(25, 9).Value2 = "PTC"

$ws.Cells.Item(26, 4).Value2 = 120.0699996948242
$ws.Cells.Item(26, 5).Value2 = 132.9100036621094
$ws.Cells.Item(26, 6).Value2 = 141.6699981689453
$ws.Cells.Item(26, 7).Value2 = 116.0199966430664
$ws.Cells.Item(26, 8).Value2 = 119792704
$ws.Cells.Item(26, 9).Value2 = "PTC"

$ws.Cells.Item(27, 4).Value2 = 139.6999969482422
$ws.Cells.Item(27, 5).Value2 = 130.9400024414062
$ws.Cells.Item(27, 6).Value2 = 149.5
$ws.Cells.Item(27, 7).Value2 = 130.5200042724609
$ws.Cells.Item(27, 8).Value2 = 119792704
$ws.Cells.Item(27, 9).Value2 = "PTC"

$ws.Cells.Item(28, 4).Value2 = 140.8800048828125
$ws.Cells.Item(28, 5).Value2 = 135.4499969482422
$ws.Cells.Item(28, 6).Value2 = 153.7299957275391
$ws.Cells.Item(28, 7).Value2 = 130.0500030517578
$ws.Cells.Item(28, 8).Value2 = 119792704
$ws.Cells.Item(28, 9).Value2 = "PTC"

$ws.Cells.Item(29, 4).Value2 = 120.4000015258789
$ws.Cells.Item(29, 5).Value2 = 127.3499984741211
$ws.Cells.Item(29, 6).Value2 = 130.3600006103516
$ws.Cells.Item(29, 7).Value2 = 114.6500015258789
$ws.Cells.Item(29, 8).Value2 = 119792704
$ws.Cells.Item(29, 9).Value2 = "PTC"

$ws.Cells.Item(30, 4).Value2 = 121.3199996948242
$ws.Cells.Item(30, 5).Value2 = 116.2600021362305
$ws.Cells.Item(30, 6).Value2 = 123.4599990844727
$ws.Cells.Item(30, 7).Value2 = 106
$ws.Cells.Item(30, 8).Value2 = 119792704
$ws.Cells.Item(30, 9).Value2 = "PTC"

$ws.Cells.Item(31, 4).Value2 = 108.5500030517578
$ws.Cells.Item(31, 5).Value2 = 114.2099990844727
$ws.Cells.Item(31, 6).Value2 = 118.4899978637695
$ws.Cells.Item(31, 7).Value2 = 96.5500030517578
$ws.Cells.Item(31, 8).Value2 = 119792704
$ws.Cells.Item(31, 9).Value2 = "PTC"

$ws.Cells.Item(32, 4).Value2 = 106.25
$ws.Cells.Item(32, 5).Value2 = 123.379997253418
$ws.Cells.Item(32, 6).Value2 = 126.879997253418
$ws.Cells.Item(32, 7).Value2 = 97.97000122070312
$ws.Cells.Item(32, 8).Value2 = 119792704
$ws.Cells.Item(32, 9).Value2 = "PTC"

$ws.Cells.Item(33, 4).Value2 = 105.8199996948242
$ws.Cells.Item(33, 5).Value2 = 117.8300018310547
$ws.Cells.Item(33, 6).Value2 = 118.8300018310547
$ws.Cells.Item(33, 7).Value2 = 101.8199996948242
$ws.Cells.Item(33, 8).Value2 = 119792704
$ws.Cells.Item(33, 9).Value2 = "PTC"

$ws.Cells.Item(34, 4).Value2 = 121.25
$ws.Cells.Item(34, 5).Value2 = 134.8800048828125
$ws.Cells.Item(34, 6).Value2 = 139.9100036621094
$ws.Cells.Item(34, 7).Value2 = 118.2099990844727
$ws.Cells.Item(34, 8).Value2 = 119792704
$ws.Cells.Item(34, 9).Value2 = "PTC"

$ws.Cells.Item(35, 4).Value2 = 127.2600021362305
$ws.Cells.Item(35, 5).Value2 = 125.7900009155273
$ws.Cells.Item(35, 6).Value2 = 130.2799987792969
$ws.Cells.Item(35, 7).Value2 = 120.620002746582
$ws.Cells.Item(35, 8).Value2 = 119792704
$ws.Cells.Item(35, 9).Value2 = "PTC"

$ws.Cells.Item(36, 4).Value2 = 142
$ws.Cells.Item(36, 5).Value2 = 145.8099975585938
$ws.Cells.Item(36, 6).Value2 = 152.0899963378906
$ws.Cells.Item(36, 7).Value2 = 135.75
$ws.Cells.Item(36, 8).Value2 = 119792704
$ws.Cells.Item(36, 9).Value2 = "PTC"

$ws.Cells.Item(37, 4).Value2 = 141.3800048828125
$ws.Cells.Item(37, 5).Value2 = 140.4199981689453
$ws.Cells.Item(37, 6).Value2 = 146.7899932861328
$ws.Cells.Item(37, 7).Value2 = 134.6100006103516
$ws.Cells.Item(37, 8).Value2 = 119792704
$ws.Cells.Item(37, 9).Value2 = "PTC"

$ws.Cells.Item(38, 4).Value2 = 173.9900054931641
$ws.Cells.Item(38, 5).Value2 = 180.6499938964844
$ws.Cells.Item(38, 6).Value2 = 185.0800018310547
$ws.Cells.Item(38, 7).Value2 = 165.7400054931641
$ws.Cells.Item(38, 8).Value2 = 119792704
$ws.Cells.Item(38, 9).Value2 = "PTC"

$ws.Cells.Item(39, 4).Value2 = 189.0899963378907
$ws.Cells.Item(39, 5).Value2 = 177.4400024414062
$ws.Cells.Item(39, 6).Value2 = 190.479995727539
$ws.Cells.Item(39, 7).Value2 = 173.8200073242188
$ws.Cells.Item(39, 8).Value2 = 119792704
$ws.Cells.Item(39, 9).Value2 = "PTC"

$ws.Cells.Item(40, 4).Value2 = 180.8300018310547
$ws.Cells.Item(40, 5).Value2 = 177.8500061035156
$ws.Cells.Item(40, 6).Value2 = 187.7799987792969
$ws.Cells.Item(40, 7).Value2 = 171.0899963378906
$ws.Cells.Item(40, 8).Value2 = 119792704
$ws.Cells.Item(40, 9).Value2 = "PTC"

$ws.Cells.Item(41, 4).Value2 = 180.0899963378906
$ws.Cells.Item(41, 5).Value2 = 185.3300018310547
$ws.Cells.Item(41, 6).Value2 = 190.8699951171875
$ws.Cells.Item(41, 7).Value2 = 174.7400054931641
$ws.Cells.Item(41, 8).Value2 = 119792704
$ws.Cells.Item(41, 9).Value2 = "PTC"

$ws.Cells.Item(42, 4).Value2 = 184.4100036621093
$ws.Cells.Item(42, 5).Value2 = 193.479995727539
$ws.Cells.Item(42, 6).Value2 = 195.229995727539
$ws.Cells.Item(42, 7).Value2 = 180.8200073242188
$ws.Cells.Item(42, 8).Value2 = 119792704
$ws.Cells.Item(42, 9).Value2 = "PTC"

$ws.Cells.Item(43, 4).Value2 = 154.1900024414062
$ws.Cells.Item(43, 5).Value2 = 154.9700012207031
$ws.Cells.Item(43, 6).Value2 = 158.4199981689453
$ws.Cells.Item(43, 7).Value2 = 133.3800048828125
$ws.Cells.Item(43, 8).Value2 = 119792704
$ws.Cells.Item(43, 9).Value2 = "PTC"

$ws.Cells.Item(44, 4).Value2 = 171.7400054931641
$ws.Cells.Item(44, 5).Value2 = 214.8099975585937
$ws.Cells.Item(44, 6).Value2 = 219.6900024414062
$ws.Cells.Item(44, 7).Value2 = 170.1399993896484
$ws.Cells.Item(44, 8).Value2 = 119792704
$ws.Cells.Item(44, 9).Value2 = "PTC"
